# Fruta / hortaliza, semanal
# The weekly price entries for rows 3-10 (excluding row 7, which is
# unchanged) were re-shuffled across their original rows. Columns A, B, C,
# E, F, G, H, I, J, K stay identical (same market/product/variety); only
# the per-record fields (Fecha, Calidad, Volumen, Precio mínimo/máximo/
# promedio ponderado, Unidad de comercialización, Origen, Precio $/Kg and
# Kg / unidad) move to different rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Record($Row, $Fecha, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $Origen, $PrecioKg, $KgUnidad) {
    $ws.Cells.Item($Row, 4).Value = $Fecha        # D: Fecha
    $ws.Cells.Item($Row, 12).Value = $Calidad     # L: Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen     # M: Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin   # N: Precio mínimo
    $ws.Cells.Item($Row, 15).Value = $PrecioMax   # O: Precio máximo
    $ws.Cells.Item($Row, 16).Value = $PrecioProm  # P: Precio promedio ponderado
    $ws.Cells.Item($Row, 17).Value = $Unidad      # Q: Unidad de comercialización
    $ws.Cells.Item($Row, 18).Value = $Origen      # R: Origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg    # S: Precio $/Kg
    $ws.Cells.Item($Row, 20).Value = $KgUnidad    # T: Kg / unidad
}

Set-Record 3 44392 "Especial" 500 7000 7000 7000 "`$/bandeja 8 kilos" "Región de O'Higgins" 875 8

Set-Record 4 44427 "Primera" 55 7000 7000 7000 "`$/caja 15 kilos granel" "Región de O'Higgins" 467 15

Set-Record 5 44495 "Primera" 50 24000 24000 24000 "`$/bandeja 10 kilos" "China" 2400 10

Set-Record 6 44418 "Especial" 100 8000 8000 8000 "`$/caja 15 kilos granel" "Región de O'Higgins" 533 15

Set-Record 8 44264 "Calibre 100" 50 20000 20000 20000 "`$/caja 18 kilos embalada" "Región de O'Higgins" 1111 18

Set-Record 9 44411 "Primera" 210 8000 8000 8000 "`$/bandeja 8 kilos" "Región de O'Higgins" 1000 8

Set-Record 10 44511 "Primera" 15 22000 22000 22000 "`$/caja 15 kilos granel" "Región de O'Higgins" 1467 15
